# Fix formatting when scraping floating point numbers:
# 1) Three "Razon social" text cells used a comma as a separator between
#    multiple names, which collided with CSV-style comma parsing downstream.
#    Replace those separating commas with periods (and normalise "S.H." -> "SH").
# 2) The "Importe" column (H) stored amounts as Spanish-formatted text
#    ("1.234,56" = thousands "." + decimal ",").  Re-write each as plain
#    decimal text ("1234.56") while keeping the cell type as text, matching
#    the sanitized OOXML in the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Razon social text fixes -----------------------------------------
$ws.Range("E44").Value  = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E89").Value  = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E177").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E174").Value = "FERNANDEZ. MARIO HUGO"
$ws.Range("E196").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"

# --- 2) Importe (column H) decimal-format fixes -------------------------
# Cells must stay TEXT (not be re-parsed as numbers), so force a text
# number format on each cell before writing the new literal string.
$c = $ws.Range("H2")
$c.NumberFormat = "@"
$c.Value = "4050.00"
$c = $ws.Range("H3")
$c.NumberFormat = "@"
$c.Value = "5360.00"
$c = $ws.Range("H4")
$c.NumberFormat = "@"
$c.Value = "1980.00"
$c = $ws.Range("H5")
$c.NumberFormat = "@"
$c.Value = "3360.00"
$c = $ws.Range("H6")
$c.NumberFormat = "@"
$c.Value = "3000.01"
$c = $ws.Range("H7")
$c.NumberFormat = "@"
$c.Value = "5300.00"
$c = $ws.Range("H8")
$c.NumberFormat = "@"
$c.Value = "200000.00"
$c = $ws.Range("H9")
$c.NumberFormat = "@"
$c.Value = "204.00"
$c = $ws.Range("H10")
$c.NumberFormat = "@"
$c.Value = "245.50"
$c = $ws.Range("H11")
$c.NumberFormat = "@"
$c.Value = "1440.00"
$c = $ws.Range("H12")
$c.NumberFormat = "@"
$c.Value = "98000.00"
$c = $ws.Range("H13")
$c.NumberFormat = "@"
$c.Value = "4660.00"
$c = $ws.Range("H14")
$c.NumberFormat = "@"
$c.Value = "9800.00"
$c = $ws.Range("H15")
$c.NumberFormat = "@"
$c.Value = "79500.00"
$c = $ws.Range("H16")
$c.NumberFormat = "@"
$c.Value = "7031.50"
$c = $ws.Range("H17")
$c.NumberFormat = "@"
$c.Value = "440638.54"
$c = $ws.Range("H18")
$c.NumberFormat = "@"
$c.Value = "248703.84"
$c = $ws.Range("H19")
$c.NumberFormat = "@"
$c.Value = "7001.50"
$c = $ws.Range("H20")
$c.NumberFormat = "@"
$c.Value = "31991.44"
$c = $ws.Range("H21")
$c.NumberFormat = "@"
$c.Value = "13120.00"
$c = $ws.Range("H22")
$c.NumberFormat = "@"
$c.Value = "559.60"
$c = $ws.Range("H23")
$c.NumberFormat = "@"
$c.Value = "10105.00"
$c = $ws.Range("H24")
$c.NumberFormat = "@"
$c.Value = "39250.00"
$c = $ws.Range("H25")
$c.NumberFormat = "@"
$c.Value = "2650.00"
$c = $ws.Range("H26")
$c.NumberFormat = "@"
$c.Value = "100.00"
$c = $ws.Range("H27")
$c.NumberFormat = "@"
$c.Value = "23411.99"
$c = $ws.Range("H28")
$c.NumberFormat = "@"
$c.Value = "19476.00"
$c = $ws.Range("H29")
$c.NumberFormat = "@"
$c.Value = "60750.00"
$c = $ws.Range("H30")
$c.NumberFormat = "@"
$c.Value = "25377.73"
$c = $ws.Range("H31")
$c.NumberFormat = "@"
$c.Value = "36209.35"
$c = $ws.Range("H32")
$c.NumberFormat = "@"
$c.Value = "4000.00"
$c = $ws.Range("H33")
$c.NumberFormat = "@"
$c.Value = "1350.00"
$c = $ws.Range("H34")
$c.NumberFormat = "@"
$c.Value = "4800.00"
$c = $ws.Range("H35")
$c.NumberFormat = "@"
$c.Value = "12000.00"
$c = $ws.Range("H36")
$c.NumberFormat = "@"
$c.Value = "144.00"
$c = $ws.Range("H37")
$c.NumberFormat = "@"
$c.Value = "400.00"
$c = $ws.Range("H38")
$c.NumberFormat = "@"
$c.Value = "947.20"
$c = $ws.Range("H39")
$c.NumberFormat = "@"
$c.Value = "168.00"
$c = $ws.Range("H40")
$c.NumberFormat = "@"
$c.Value = "40.00"
$c = $ws.Range("H41")
$c.NumberFormat = "@"
$c.Value = "52670.18"
$c = $ws.Range("H42")
$c.NumberFormat = "@"
$c.Value = "22451.25"
$c = $ws.Range("H43")
$c.NumberFormat = "@"
$c.Value = "331.24"
$c = $ws.Range("H44")
$c.NumberFormat = "@"
$c.Value = "1910.00"
$c = $ws.Range("H45")
$c.NumberFormat = "@"
$c.Value = "1130.85"
$c = $ws.Range("H46")
$c.NumberFormat = "@"
$c.Value = "5375.00"
$c = $ws.Range("H47")
$c.NumberFormat = "@"
$c.Value = "26100.00"
$c = $ws.Range("H48")
$c.NumberFormat = "@"
$c.Value = "73620.70"
$c = $ws.Range("H49")
$c.NumberFormat = "@"
$c.Value = "4287.97"
$c = $ws.Range("H50")
$c.NumberFormat = "@"
$c.Value = "790.00"
$c = $ws.Range("H51")
$c.NumberFormat = "@"
$c.Value = "375.00"
$c = $ws.Range("H52")
$c.NumberFormat = "@"
$c.Value = "75870.93"
$c = $ws.Range("H53")
$c.NumberFormat = "@"
$c.Value = "11757.60"
$c = $ws.Range("H54")
$c.NumberFormat = "@"
$c.Value = "186.00"
$c = $ws.Range("H55")
$c.NumberFormat = "@"
$c.Value = "4788.20"
$c = $ws.Range("H56")
$c.NumberFormat = "@"
$c.Value = "5579.00"
$c = $ws.Range("H57")
$c.NumberFormat = "@"
$c.Value = "7167.54"
$c = $ws.Range("H58")
$c.NumberFormat = "@"
$c.Value = "732.26"
$c = $ws.Range("H59")
$c.NumberFormat = "@"
$c.Value = "44398.00"
$c = $ws.Range("H60")
$c.NumberFormat = "@"
$c.Value = "24.06"
$c = $ws.Range("H61")
$c.NumberFormat = "@"
$c.Value = "343200.00"
$c = $ws.Range("H62")
$c.NumberFormat = "@"
$c.Value = "4701.80"
$c = $ws.Range("H63")
$c.NumberFormat = "@"
$c.Value = "3201.18"
$c = $ws.Range("H64")
$c.NumberFormat = "@"
$c.Value = "5500.00"
$c = $ws.Range("H65")
$c.NumberFormat = "@"
$c.Value = "3374.00"
$c = $ws.Range("H66")
$c.NumberFormat = "@"
$c.Value = "6049.00"
$c = $ws.Range("H67")
$c.NumberFormat = "@"
$c.Value = "5680.51"
$c = $ws.Range("H68")
$c.NumberFormat = "@"
$c.Value = "2436.53"
$c = $ws.Range("H69")
$c.NumberFormat = "@"
$c.Value = "1310.00"
$c = $ws.Range("H70")
$c.NumberFormat = "@"
$c.Value = "780.00"
$c = $ws.Range("H71")
$c.NumberFormat = "@"
$c.Value = "1500.00"
$c = $ws.Range("H72")
$c.NumberFormat = "@"
$c.Value = "22650.00"
$c = $ws.Range("H73")
$c.NumberFormat = "@"
$c.Value = "5534.44"
$c = $ws.Range("H74")
$c.NumberFormat = "@"
$c.Value = "1432.00"
$c = $ws.Range("H75")
$c.NumberFormat = "@"
$c.Value = "5280.00"
$c = $ws.Range("H76")
$c.NumberFormat = "@"
$c.Value = "48.00"
$c = $ws.Range("H77")
$c.NumberFormat = "@"
$c.Value = "5874.33"
$c = $ws.Range("H78")
$c.NumberFormat = "@"
$c.Value = "54410.00"
$c = $ws.Range("H79")
$c.NumberFormat = "@"
$c.Value = "2570.00"
$c = $ws.Range("H80")
$c.NumberFormat = "@"
$c.Value = "14999.85"
$c = $ws.Range("H81")
$c.NumberFormat = "@"
$c.Value = "848.00"
$c = $ws.Range("H82")
$c.NumberFormat = "@"
$c.Value = "21000.00"
$c = $ws.Range("H83")
$c.NumberFormat = "@"
$c.Value = "52540.00"
$c = $ws.Range("H84")
$c.NumberFormat = "@"
$c.Value = "25128.32"
$c = $ws.Range("H85")
$c.NumberFormat = "@"
$c.Value = "4.12"
$c = $ws.Range("H86")
$c.NumberFormat = "@"
$c.Value = "9690.00"
$c = $ws.Range("H87")
$c.NumberFormat = "@"
$c.Value = "385.00"
$c = $ws.Range("H88")
$c.NumberFormat = "@"
$c.Value = "1200.00"
$c = $ws.Range("H89")
$c.NumberFormat = "@"
$c.Value = "4360.00"
$c = $ws.Range("H90")
$c.NumberFormat = "@"
$c.Value = "19045.00"
$c = $ws.Range("H91")
$c.NumberFormat = "@"
$c.Value = "172.80"
$c = $ws.Range("H92")
$c.NumberFormat = "@"
$c.Value = "2150.00"
$c = $ws.Range("H93")
$c.NumberFormat = "@"
$c.Value = "14908.98"
$c = $ws.Range("H94")
$c.NumberFormat = "@"
$c.Value = "5289.00"
$c = $ws.Range("H95")
$c.NumberFormat = "@"
$c.Value = "351.48"
$c = $ws.Range("H96")
$c.NumberFormat = "@"
$c.Value = "12053.82"
$c = $ws.Range("H97")
$c.NumberFormat = "@"
$c.Value = "3.54"
$c = $ws.Range("H98")
$c.NumberFormat = "@"
$c.Value = "526097.59"
$c = $ws.Range("H99")
$c.NumberFormat = "@"
$c.Value = "53261.49"
$c = $ws.Range("H100")
$c.NumberFormat = "@"
$c.Value = "39.69"
$c = $ws.Range("H101")
$c.NumberFormat = "@"
$c.Value = "80.32"
$c = $ws.Range("H102")
$c.NumberFormat = "@"
$c.Value = "29752.96"
$c = $ws.Range("H103")
$c.NumberFormat = "@"
$c.Value = "288.00"
$c = $ws.Range("H104")
$c.NumberFormat = "@"
$c.Value = "7339.00"
$c = $ws.Range("H105")
$c.NumberFormat = "@"
$c.Value = "45.00"
$c = $ws.Range("H106")
$c.NumberFormat = "@"
$c.Value = "1567.23"
$c = $ws.Range("H107")
$c.NumberFormat = "@"
$c.Value = "2929.87"
$c = $ws.Range("H108")
$c.NumberFormat = "@"
$c.Value = "3912.60"
$c = $ws.Range("H109")
$c.NumberFormat = "@"
$c.Value = "4760.00"
$c = $ws.Range("H110")
$c.NumberFormat = "@"
$c.Value = "18971.98"
$c = $ws.Range("H111")
$c.NumberFormat = "@"
$c.Value = "2960.00"
$c = $ws.Range("H112")
$c.NumberFormat = "@"
$c.Value = "431.22"
$c = $ws.Range("H113")
$c.NumberFormat = "@"
$c.Value = "7070.00"
$c = $ws.Range("H114")
$c.NumberFormat = "@"
$c.Value = "2240.00"
$c = $ws.Range("H115")
$c.NumberFormat = "@"
$c.Value = "2682.00"
$c = $ws.Range("H116")
$c.NumberFormat = "@"
$c.Value = "390.00"
$c = $ws.Range("H117")
$c.NumberFormat = "@"
$c.Value = "32751.00"
$c = $ws.Range("H118")
$c.NumberFormat = "@"
$c.Value = "69.00"
$c = $ws.Range("H119")
$c.NumberFormat = "@"
$c.Value = "16053.06"
$c = $ws.Range("H120")
$c.NumberFormat = "@"
$c.Value = "6106.00"
$c = $ws.Range("H121")
$c.NumberFormat = "@"
$c.Value = "2761.14"
$c = $ws.Range("H122")
$c.NumberFormat = "@"
$c.Value = "2639.20"
$c = $ws.Range("H123")
$c.NumberFormat = "@"
$c.Value = "771.05"
$c = $ws.Range("H124")
$c.NumberFormat = "@"
$c.Value = "324.00"
$c = $ws.Range("H125")
$c.NumberFormat = "@"
$c.Value = "321.00"
$c = $ws.Range("H126")
$c.NumberFormat = "@"
$c.Value = "134.10"
$c = $ws.Range("H127")
$c.NumberFormat = "@"
$c.Value = "7100.00"
$c = $ws.Range("H128")
$c.NumberFormat = "@"
$c.Value = "6900.00"
$c = $ws.Range("H129")
$c.NumberFormat = "@"
$c.Value = "1008.70"
$c = $ws.Range("H130")
$c.NumberFormat = "@"
$c.Value = "500.00"
$c = $ws.Range("H131")
$c.NumberFormat = "@"
$c.Value = "237.00"
$c = $ws.Range("H132")
$c.NumberFormat = "@"
$c.Value = "255.00"
$c = $ws.Range("H133")
$c.NumberFormat = "@"
$c.Value = "2400.00"
$c = $ws.Range("H134")
$c.NumberFormat = "@"
$c.Value = "2300.00"
$c = $ws.Range("H135")
$c.NumberFormat = "@"
$c.Value = "34700.00"
$c = $ws.Range("H136")
$c.NumberFormat = "@"
$c.Value = "6486.00"
$c = $ws.Range("H137")
$c.NumberFormat = "@"
$c.Value = "5483.00"
$c = $ws.Range("H138")
$c.NumberFormat = "@"
$c.Value = "6910.00"
$c = $ws.Range("H139")
$c.NumberFormat = "@"
$c.Value = "2328.00"
$c = $ws.Range("H140")
$c.NumberFormat = "@"
$c.Value = "6000.00"
$c = $ws.Range("H141")
$c.NumberFormat = "@"
$c.Value = "28500.00"
$c = $ws.Range("H142")
$c.NumberFormat = "@"
$c.Value = "8500.00"
$c = $ws.Range("H143")
$c.NumberFormat = "@"
$c.Value = "114100.00"
$c = $ws.Range("H144")
$c.NumberFormat = "@"
$c.Value = "4500.00"
$c = $ws.Range("H145")
$c.NumberFormat = "@"
$c.Value = "1917.61"
$c = $ws.Range("H146")
$c.NumberFormat = "@"
$c.Value = "240.68"
$c = $ws.Range("H147")
$c.NumberFormat = "@"
$c.Value = "6076.00"
$c = $ws.Range("H148")
$c.NumberFormat = "@"
$c.Value = "2044680.00"
$c = $ws.Range("H149")
$c.NumberFormat = "@"
$c.Value = "35943.70"
$c = $ws.Range("H150")
$c.NumberFormat = "@"
$c.Value = "8000.00"
$c = $ws.Range("H151")
$c.NumberFormat = "@"
$c.Value = "3500.00"
$c = $ws.Range("H152")
$c.NumberFormat = "@"
$c.Value = "10000.00"
$c = $ws.Range("H153")
$c.NumberFormat = "@"
$c.Value = "32505.44"
$c = $ws.Range("H154")
$c.NumberFormat = "@"
$c.Value = "2500.00"
$c = $ws.Range("H155")
$c.NumberFormat = "@"
$c.Value = "3204.50"
$c = $ws.Range("H156")
$c.NumberFormat = "@"
$c.Value = "2556.00"
$c = $ws.Range("H157")
$c.NumberFormat = "@"
$c.Value = "3000.00"
$c = $ws.Range("H158")
$c.NumberFormat = "@"
$c.Value = "3000.00"
$c = $ws.Range("H159")
$c.NumberFormat = "@"
$c.Value = "2000.00"
$c = $ws.Range("H160")
$c.NumberFormat = "@"
$c.Value = "1500.00"
$c = $ws.Range("H161")
$c.NumberFormat = "@"
$c.Value = "17874.50"
$c = $ws.Range("H162")
$c.NumberFormat = "@"
$c.Value = "7500.00"
$c = $ws.Range("H163")
$c.NumberFormat = "@"
$c.Value = "4000.00"
$c = $ws.Range("H164")
$c.NumberFormat = "@"
$c.Value = "2500.00"
$c = $ws.Range("H165")
$c.NumberFormat = "@"
$c.Value = "4100.00"
$c = $ws.Range("H166")
$c.NumberFormat = "@"
$c.Value = "19810.53"
$c = $ws.Range("H167")
$c.NumberFormat = "@"
$c.Value = "6570.00"
$c = $ws.Range("H168")
$c.NumberFormat = "@"
$c.Value = "6000.00"
$c = $ws.Range("H169")
$c.NumberFormat = "@"
$c.Value = "2500.00"
$c = $ws.Range("H170")
$c.NumberFormat = "@"
$c.Value = "32280.00"
$c = $ws.Range("H171")
$c.NumberFormat = "@"
$c.Value = "4000.00"
$c = $ws.Range("H172")
$c.NumberFormat = "@"
$c.Value = "680.00"
$c = $ws.Range("H173")
$c.NumberFormat = "@"
$c.Value = "585.00"
$c = $ws.Range("H174")
$c.NumberFormat = "@"
$c.Value = "380.00"
$c = $ws.Range("H175")
$c.NumberFormat = "@"
$c.Value = "7440.00"
$c = $ws.Range("H176")
$c.NumberFormat = "@"
$c.Value = "57218.80"
$c = $ws.Range("H177")
$c.NumberFormat = "@"
$c.Value = "700.00"
$c = $ws.Range("H178")
$c.NumberFormat = "@"
$c.Value = "740.00"
$c = $ws.Range("H179")
$c.NumberFormat = "@"
$c.Value = "2953.00"
$c = $ws.Range("H180")
$c.NumberFormat = "@"
$c.Value = "1900.00"
$c = $ws.Range("H181")
$c.NumberFormat = "@"
$c.Value = "1025.00"
$c = $ws.Range("H182")
$c.NumberFormat = "@"
$c.Value = "85800.00"
$c = $ws.Range("H183")
$c.NumberFormat = "@"
$c.Value = "650.00"
$c = $ws.Range("H184")
$c.NumberFormat = "@"
$c.Value = "1521.30"
$c = $ws.Range("H185")
$c.NumberFormat = "@"
$c.Value = "14600.00"
$c = $ws.Range("H186")
$c.NumberFormat = "@"
$c.Value = "12950.00"
$c = $ws.Range("H187")
$c.NumberFormat = "@"
$c.Value = "570.00"
$c = $ws.Range("H188")
$c.NumberFormat = "@"
$c.Value = "67.29"
$c = $ws.Range("H189")
$c.NumberFormat = "@"
$c.Value = "4670.00"
$c = $ws.Range("H190")
$c.NumberFormat = "@"
$c.Value = "9500.00"
$c = $ws.Range("H191")
$c.NumberFormat = "@"
$c.Value = "2422.00"
$c = $ws.Range("H192")
$c.NumberFormat = "@"
$c.Value = "3453.00"
$c = $ws.Range("H193")
$c.NumberFormat = "@"
$c.Value = "6791.00"
$c = $ws.Range("H194")
$c.NumberFormat = "@"
$c.Value = "740.00"
$c = $ws.Range("H195")
$c.NumberFormat = "@"
$c.Value = "12133.00"
$c = $ws.Range("H196")
$c.NumberFormat = "@"
$c.Value = "7300.00"
$c = $ws.Range("H197")
$c.NumberFormat = "@"
$c.Value = "1250.00"
$c = $ws.Range("H198")
$c.NumberFormat = "@"
$c.Value = "4471.60"
$c = $ws.Range("H199")
$c.NumberFormat = "@"
$c.Value = "428.16"
$c = $ws.Range("H200")
$c.NumberFormat = "@"
$c.Value = "67743.12"
$c = $ws.Range("H201")
$c.NumberFormat = "@"
$c.Value = "359.16"
$c = $ws.Range("H202")
$c.NumberFormat = "@"
$c.Value = "3760.00"
$c = $ws.Range("H203")
$c.NumberFormat = "@"
$c.Value = "300.00"
$c = $ws.Range("H204")
$c.NumberFormat = "@"
$c.Value = "122826.00"
$c = $ws.Range("H205")
$c.NumberFormat = "@"
$c.Value = "800.00"
$c = $ws.Range("H206")
$c.NumberFormat = "@"
$c.Value = "12500.00"
$c = $ws.Range("H207")
$c.NumberFormat = "@"
$c.Value = "44301.25"
$c = $ws.Range("H208")
$c.NumberFormat = "@"
$c.Value = "25000.00"
$c = $ws.Range("H209")
$c.NumberFormat = "@"
$c.Value = "25000.00"
$c = $ws.Range("H210")
$c.NumberFormat = "@"
$c.Value = "25000.00"
$c = $ws.Range("H211")
$c.NumberFormat = "@"
$c.Value = "25000.00"
$c = $ws.Range("H212")
$c.NumberFormat = "@"
$c.Value = "50000.00"
$c = $ws.Range("H213")
$c.NumberFormat = "@"
$c.Value = "50000.00"
$c = $ws.Range("H214")
$c.NumberFormat = "@"
$c.Value = "25000.00"
$c = $ws.Range("H215")
$c.NumberFormat = "@"
$c.Value = "42000.00"
$c = $ws.Range("H216")
$c.NumberFormat = "@"
$c.Value = "9138.38"
$c = $ws.Range("H217")
$c.NumberFormat = "@"
$c.Value = "4140213.30"
$c = $ws.Range("H218")
$c.NumberFormat = "@"
$c.Value = "10892.91"
$c = $ws.Range("H219")
$c.NumberFormat = "@"
$c.Value = "1800.00"
$c = $ws.Range("H220")
$c.NumberFormat = "@"
$c.Value = "83100.00"
$c = $ws.Range("H221")
$c.NumberFormat = "@"
$c.Value = "230000.00"
$c = $ws.Range("H222")
$c.NumberFormat = "@"
$c.Value = "105000.00"
$c = $ws.Range("H223")
$c.NumberFormat = "@"
$c.Value = "105000.00"
$c = $ws.Range("H224")
$c.NumberFormat = "@"
$c.Value = "105000.00"
$c = $ws.Range("H225")
$c.NumberFormat = "@"
$c.Value = "105000.00"
$c = $ws.Range("H226")
$c.NumberFormat = "@"
$c.Value = "105000.00"
$c = $ws.Range("H227")
$c.NumberFormat = "@"
$c.Value = "175000.00"
$c = $ws.Range("H228")
$c.NumberFormat = "@"
$c.Value = "175000.00"
$c = $ws.Range("H229")
$c.NumberFormat = "@"
$c.Value = "245000.00"
$c = $ws.Range("H230")
$c.NumberFormat = "@"
$c.Value = "105000.00"
$c = $ws.Range("H231")
$c.NumberFormat = "@"
$c.Value = "105000.00"
$c = $ws.Range("H232")
$c.NumberFormat = "@"
$c.Value = "105000.00"
$c = $ws.Range("H233")
$c.NumberFormat = "@"
$c.Value = "105000.00"
$c = $ws.Range("H234")
$c.NumberFormat = "@"
$c.Value = "105000.00"
$c = $ws.Range("H235")
$c.NumberFormat = "@"
$c.Value = "175000.00"
$c = $ws.Range("H236")
$c.NumberFormat = "@"
$c.Value = "315000.00"
$c = $ws.Range("H237")
$c.NumberFormat = "@"
$c.Value = "175000.00"
$c = $ws.Range("H238")
$c.NumberFormat = "@"
$c.Value = "105000.00"
$c = $ws.Range("H239")
$c.NumberFormat = "@"
$c.Value = "155000.00"
$c = $ws.Range("H240")
$c.NumberFormat = "@"
$c.Value = "105000.00"
$c = $ws.Range("H241")
$c.NumberFormat = "@"
$c.Value = "105000.00"
$c = $ws.Range("H242")
$c.NumberFormat = "@"
$c.Value = "105000.00"
$c = $ws.Range("H243")
$c.NumberFormat = "@"
$c.Value = "105000.00"
$c = $ws.Range("H244")
$c.NumberFormat = "@"
$c.Value = "70950.35"
$c = $ws.Range("H245")
$c.NumberFormat = "@"
$c.Value = "20600.00"
$c = $ws.Range("H246")
$c.NumberFormat = "@"
$c.Value = "5000.00"
$c = $ws.Range("H247")
$c.NumberFormat = "@"
$c.Value = "78287.00"
$c = $ws.Range("H248")
$c.NumberFormat = "@"
$c.Value = "6590.00"
$c = $ws.Range("H249")
$c.NumberFormat = "@"
$c.Value = "7000.00"
$c = $ws.Range("H250")
$c.NumberFormat = "@"
$c.Value = "415.63"
$c = $ws.Range("H251")
$c.NumberFormat = "@"
$c.Value = "3000.00"
$c = $ws.Range("H252")
$c.NumberFormat = "@"
$c.Value = "278633.16"
$c = $ws.Range("H253")
$c.NumberFormat = "@"
$c.Value = "15085.00"
$c = $ws.Range("H254")
$c.NumberFormat = "@"
$c.Value = "7000.02"
$c = $ws.Range("H255")
$c.NumberFormat = "@"
$c.Value = "34807.50"
$c = $ws.Range("H256")
$c.NumberFormat = "@"
$c.Value = "30000.00"
$c = $ws.Range("H257")
$c.NumberFormat = "@"
$c.Value = "6000.00"
$c = $ws.Range("H258")
$c.NumberFormat = "@"
$c.Value = "36300.00"
$c = $ws.Range("H259")
$c.NumberFormat = "@"
$c.Value = "1110.00"
